# syllabus-2016.xlsx update: "More reorganization, some updates to python 3"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Window position (cosmetic; best-effort) ---
try {
    $wb.Windows.Item(1).Left = 3400
    $wb.Windows.Item(1).Top = 2500
} catch {}

# --- Cell content updates -------------------------------------------------

# Week 1 (row 3) / Week 2 (rows 4-5)
$ws.Range("G3").Value = "Introduction-to-jupyter-notebooks-and-python.ipynb"
$ws.Range("G4").Value = "Strings-lists-dictionaries.ipynb"
$ws.Range("D5").Value = "Wordpress tutorial"

# Week 3 (row 7) -> add notebook reference in G7, matching style of neighboring cells
$ws.Range("F7").Copy()
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("G7").Value = "?"

# Week 4 (rows 8-9)
$ws.Range("D8").Value = "String Methods, Conditions, Booleans and Iteration"
$ws.Range("F8").Copy()
$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("G8").Value = "string-methods-conditions-iteration.ipynb"

$ws.Range("E9").Value = "[PW] ET Confirmed"
$ws.Range("F9").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("G9").Value = "NA"

# Week 5 (row 10)
$ws.Range("D10").Value = "From  Lists, Data Cleaning"
$ws.Range("G10").Value = "Planning-code-prime-numbers-example.ipynb"

# Week 7 (rows 15-16)
$ws.Range("E15").Value = "[SM] assign exercise 6: Extract data via API"
$ws.Range("D16").Value = "Web Programming with Javascript"
$ws.Range("E16").Value = "[GB]"

# Week 9 (row 19)
$ws.Range("D19").Value = "Web Mapping with Leaflet"
$ws.Range("E19").Value = "[GB] assign exercise 8: web mapping blog post"

# Week 11/12 reorg (rows 25-27)
$ws.Range("D26").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("D25").Value = "Geocoding and Reverse Geocoding"
$ws.Range("E25").Value = "[GB]"

$ws.Range("D16").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("D26").Value = "Github and Version Control"
$ws.Range("E26").Value = "[PW]"

$ws.Range("D16").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = "Statistical Analysis with Statsmodels"
$ws.Range("E27").Value = "[PW]"

# Week 13 (row 28)
$ws.Range("E28").Value = "[PW/SB]"

# --- Remove vestigial empty row 37 ---
$ws.Range("D37").ClearContents()

# --- Selection ---
$ws.Range("E28").Select()
